$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/A008_B009_1025IF_001.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/A008_B009_1025IF_001.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/d2 copy/efgh.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/d2 copy/efgh.mp4"),
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/A008_B009_1025IF_001.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/A008_B009_1025IF_001.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/d2 copy/efgh.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/d2 copy/efgh.mp4"),
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/A008_B009_1025IF_001.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/A008_B009_1025IF_001.mp4"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/d2 copy/efgh.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/d2 copy/efgh.mp4"),
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mov"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/A008_B009_1025IF_001.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/A008_B009_1025IF_001.mov"),
    @("/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/d2 copy/efgh.R3D", "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/d2 copy/efgh.mov"),
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mov"),
    @("/Users/rayan/Downloads/d1/abcd.R3D", "/Users/rayan/Downloads/d1_converted/abcd.mov")
)

$startRow = 8
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}
